$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 2513
$ws.Range("I61").Value = 2513
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 7539
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -7367
$ws.Range("N61").ClearContents()

$ws.Range("H86").Value = 19247.268
$ws.Range("I86").Value = 5753.5835
$ws.Range("K86").Value = 5753.5835
$ws.Range("M86").Value = -4630.5835

$ws.Range("H89").Value = 19247.268
$ws.Range("I89").Value = 5753.5835
$ws.Range("K89").Value = 28767.9175
$ws.Range("M89").Value = -23151.9175

$ws.Range("H113").Value = 8250.375
$ws.Range("I113").Value = 12626
$ws.Range("J113").Value = 3874.75
$ws.Range("K113").Value = 12626
$ws.Range("L113").Value = 3874.75
$ws.Range("M113").Value = -9372
$ws.Range("N113").Value = -10382.75

$ws.Range("H125").Value = 2529.8572
$ws.Range("I125").Value = 827.6667
$ws.Range("J125").Value = 3806.5
$ws.Range("K125").Value = 7449.0003
$ws.Range("L125").Value = 34258.5
$ws.Range("M125").Value = -4989.0003
$ws.Range("N125").Value = -39178.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3923.535
$ws.Range("I32").Value = 3494.3333
$ws.Range("K32").Value = 3494.3333
$ws.Range("M32").Value = -3207.3333

$ws.Range("H45").Value = 11318.23
$ws.Range("I45").Value = 27116.4
$ws.Range("J45").Value = 1444.375
$ws.Range("K45").Value = 27116.4
$ws.Range("L45").Value = 1444.375
$ws.Range("M45").Value = -26739.4
$ws.Range("N45").Value = -2198.375

$ws.Range("H69").Value = 179977
$ws.Range("J69").Value = 179977
$ws.Range("L69").Value = 179977
$ws.Range("N69").Value = -181475

$ws.Range("H72").Value = 179977
$ws.Range("J72").Value = 179977
$ws.Range("L72").Value = 539931
$ws.Range("N72").Value = -547419

$ws.Range("H74").Value = 7264.875
$ws.Range("I74").Value = 8032
$ws.Range("J74").Value = 5577.2
$ws.Range("K74").Value = 8032
$ws.Range("L74").Value = 5577.2
$ws.Range("M74").Value = -7158
$ws.Range("N74").Value = -7325.2

$ws.Range("H77").Value = 7264.875
$ws.Range("I77").Value = 8032
$ws.Range("J77").Value = 5577.2
$ws.Range("K77").Value = 40160
$ws.Range("L77").Value = 27886
$ws.Range("M77").Value = -35792
$ws.Range("N77").Value = -36622

$ws.Range("H132").Value = 2018.3903
$ws.Range("I132").Value = 2005.7931
$ws.Range("J132").Value = 2048.8333
$ws.Range("K132").Value = 6017.379300000001
$ws.Range("L132").Value = 6146.499899999999
$ws.Range("M132").Value = -3487.379300000001
$ws.Range("N132").Value = -11206.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H41").Value = 175677
$ws.Range("J41").Value = 175677
$ws.Range("L41").Value = 175677
$ws.Range("N41").Value = -176453

$ws.Range("H42").Value = 99999
$ws.Range("J42").Value = 99999
$ws.Range("L42").Value = 99999
$ws.Range("N42").Value = -100655

$ws.Range("H86").Value = 1798.2273
$ws.Range("I86").Value = 1691.1875
$ws.Range("J86").Value = 2083.6667
$ws.Range("K86").Value = 1691.1875
$ws.Range("L86").Value = 2083.6667
$ws.Range("M86").Value = -568.1875
$ws.Range("N86").Value = -4329.6667

$ws.Range("H89").Value = 1798.2273
$ws.Range("I89").Value = 1691.1875
$ws.Range("J89").Value = 2083.6667
$ws.Range("K89").Value = 8455.9375
$ws.Range("L89").Value = 10418.3335
$ws.Range("M89").Value = -2839.9375
$ws.Range("N89").Value = -21650.3335

$ws.Range("H94").Value = 1555.3636
$ws.Range("I94").Value = 699.8889
$ws.Range("J94").Value = 5405
$ws.Range("K94").Value = 699.8889
$ws.Range("L94").Value = 5405
$ws.Range("M94").Value = -248.8889
$ws.Range("N94").Value = -6307

$ws.Range("H107").Value = 3836.15
$ws.Range("I107").Value = 4500.8335
$ws.Range("J107").Value = 2839.125
$ws.Range("K107").Value = 4500.8335
$ws.Range("L107").Value = 2839.125
$ws.Range("M107").Value = -2580.8335
$ws.Range("N107").Value = -6679.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2799.257
$ws.Range("I31").Value = 2770.0833
$ws.Range("K31").Value = 2770.0833
$ws.Range("M31").Value = -2475.0833

$ws.Range("H34").Value = 2799.257
$ws.Range("I34").Value = 2770.0833
$ws.Range("K34").Value = 2770.0833
$ws.Range("M34").Value = -2568.0833

$ws.Range("H107").Value = 965.05884
$ws.Range("I107").Value = 678.1875
$ws.Range("K107").Value = 678.1875
$ws.Range("M107").Value = 1241.8125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 97222330
$ws.Range("I2").Value = 194444510
$ws.Range("J2").Value = 133.5
$ws.Range("K2").Value = 194444510
$ws.Range("L2").Value = 133.5
$ws.Range("M2").Value = -194444397
$ws.Range("N2").Value = -359.5

$ws.Range("H70").Value = 2359.4443
$ws.Range("I70").Value = 2470.138
$ws.Range("J70").Value = 1900.8572
$ws.Range("K70").Value = 2470.138
$ws.Range("L70").Value = 1900.8572
$ws.Range("M70").Value = -2200.138
$ws.Range("N70").Value = -2440.8572

$ws.Range("H73").Value = 2359.4443
$ws.Range("I73").Value = 2470.138
$ws.Range("J73").Value = 1900.8572
$ws.Range("K73").Value = 2470.138
$ws.Range("L73").Value = 1900.8572
$ws.Range("M73").Value = -1534.138
$ws.Range("N73").Value = -3772.8572

$ws.Range("H97").Value = 861.8148
$ws.Range("I97").Value = 391
$ws.Range("K97").Value = 391
$ws.Range("M97").Value = 105

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1745.0571
$ws.Range("I16").Value = 1455.4062
$ws.Range("J16").Value = 4834.6665
$ws.Range("K16").Value = 1455.4062
$ws.Range("L16").Value = 4834.6665
$ws.Range("M16").Value = -1285.4062
$ws.Range("N16").Value = -5174.6665

$ws.Range("H43").Value = 16999.5
$ws.Range("J43").Value = 16999.5
$ws.Range("L43").Value = 16999.5
$ws.Range("N43").Value = -17385.5

$ws.Range("H45").Value = 19397
$ws.Range("I45").Value = 17996.25
$ws.Range("K45").Value = 17996.25
$ws.Range("M45").Value = -17589.25

$ws.Range("H93").Value = 4159.3335
$ws.Range("I93").Value = 4159.3335
$ws.Range("K93").Value = 4159.3335
$ws.Range("M93").Value = -2911.3335

$ws.Range("H132").Value = 40737.55
$ws.Range("J132").Value = 3966
$ws.Range("L132").Value = 11898
$ws.Range("N132").Value = -16958

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 333338000
$ws.Range("I62").Value = 500004500
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 500004500
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -500003876
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 333338000
$ws.Range("I65").Value = 500004500
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 2500022500
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -2500019380
$ws.Range("N65").Value = -31240

$ws.Range("H113").Value = 626
$ws.Range("I113").Value = 365.2
$ws.Range("J113").Value = 1060.6666
$ws.Range("K113").Value = 1095.6
$ws.Range("L113").Value = 3181.9998
$ws.Range("M113").Value = 1074.4
$ws.Range("N113").Value = -7521.9998

$ws.Range("H126").Value = 2159234.8
$ws.Range("I126").Value = 2607970.5
$ws.Range("J126").Value = 5303.7
$ws.Range("K126").Value = 7823911.5
$ws.Range("L126").Value = 15911.1
$ws.Range("M126").Value = -7821441.5
$ws.Range("N126").Value = -20851.1

Write-Host "Applied 34 row updates across ALC, ARM, BSM, CRP, GSM, LTW, WVR sheets"
